$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$ws.Cells.Item(2, 4).Value = "27.441.76"
$ws.Cells.Item(2, 5).Value = "  +2.01%  "
$ws.Cells.Item(3, 4).Value = "1.837.41"
$ws.Cells.Item(3, 5).Value = "  +1.26%  "
Set-TextValue $ws.Cells.Item(4, 4) "1.013"
$ws.Cells.Item(4, 5).Value = "  +1.11%  "
Set-TextValue $ws.Cells.Item(5, 4) "314.29"
$ws.Cells.Item(5, 5).Value = "  +1.61%  "
$ws.Cells.Item(6, 5).Value = "  +0.89%  "
Set-TextValue $ws.Cells.Item(7, 4) "0.4740"
$ws.Cells.Item(7, 5).Value = "  +1.68%  "
Set-TextValue $ws.Cells.Item(8, 4) "0.3689"
$ws.Cells.Item(8, 5).Value = "  +0.76%  "
Set-TextValue $ws.Cells.Item(9, 4) "0.07457"
$ws.Cells.Item(9, 5).Value = "  +1.50%  "
Set-TextValue $ws.Cells.Item(10, 4) "0.8849"
$ws.Cells.Item(10, 5).Value = "  +1.85%  "
Set-TextValue $ws.Cells.Item(11, 4) "20.44"
$ws.Cells.Item(11, 5).Value = "  +0.71%  "
$ws.Cells.Item(12, 4).Value = "1.937.51"
$ws.Cells.Item(12, 5).Value = "  +4.28%  "
$ws.Cells.Item(13, 5).Value = "  +3.42%  "
Set-TextValue $ws.Cells.Item(14, 4) "5.449"
$ws.Cells.Item(14, 5).Value = "  +1.23%  "
Set-TextValue $ws.Cells.Item(15, 4) "93.23"
$ws.Cells.Item(15, 5).Value = "  +1.75%  "
Set-TextValue $ws.Cells.Item(16, 4) "6.577"
$ws.Cells.Item(16, 5).Value = "  +0.91%  "
$ws.Cells.Item(17, 5).Value = "  +0.73%  "
Set-TextValue $ws.Cells.Item(18, 4) "0.000008811"
$ws.Cells.Item(18, 5).Value = "  +1.17%  "
Set-TextValue $ws.Cells.Item(19, 4) "1.011"
$ws.Cells.Item(19, 5).Value = "  +0.93%  "
Set-TextValue $ws.Cells.Item(20, 4) "14.80"
$ws.Cells.Item(20, 5).Value = "  +0.97%  "
$ws.Cells.Item(21, 4).Value = "27.478.26"
$ws.Cells.Item(21, 5).Value = "  +2.00%  "
Set-TextValue $ws.Cells.Item(22, 4) "5.325"
$ws.Cells.Item(22, 5).Value = "  +0.44%  "
Set-TextValue $ws.Cells.Item(23, 4) "10.69"
$ws.Cells.Item(23, 5).Value = "  +0.54%  "
$ws.Cells.Item(24, 4).Value = "2.156.69"
$ws.Cells.Item(24, 5).Value = "  +4.57%  "
Set-TextValue $ws.Cells.Item(25, 4) "1.913"
$ws.Cells.Item(25, 5).Value = "  +0.94%  "
Set-TextValue $ws.Cells.Item(26, 4) "152.07"
$ws.Cells.Item(26, 5).Value = "  +0.75%  "
Set-TextValue $ws.Cells.Item(27, 4) "18.65"
$ws.Cells.Item(27, 5).Value = "  +1.96%  "
Set-TextValue $ws.Cells.Item(28, 4) "2.145"
$ws.Cells.Item(28, 5).Value = "  -0.33%  "
Set-TextValue $ws.Cells.Item(29, 4) "5.247"
$ws.Cells.Item(29, 5).Value = "  -0.29%  "
Set-TextValue $ws.Cells.Item(30, 4) "117.88"
$ws.Cells.Item(30, 5).Value = "  +2.11%  "
Set-TextValue $ws.Cells.Item(31, 4) "0.08998"
$ws.Cells.Item(31, 5).Value = "  +0.90%  "
Set-TextValue $ws.Cells.Item(32, 4) "0.7576"
$ws.Cells.Item(32, 5).Value = "  +0.43%  "
Set-TextValue $ws.Cells.Item(33, 4) "1.182"
$ws.Cells.Item(33, 5).Value = "  +2.21%  "
Set-TextValue $ws.Cells.Item(34, 4) "4.558"
$ws.Cells.Item(34, 5).Value = "  +1.59%  "
Set-TextValue $ws.Cells.Item(35, 4) "2.949"
$ws.Cells.Item(35, 5).Value = "  +1.24%  "
Set-TextValue $ws.Cells.Item(36, 4) "1.012"
$ws.Cells.Item(36, 5).Value = "  +1.07%  "
$ws.Cells.Item(37, 5).Value = "  +2.01%  "
Set-TextValue $ws.Cells.Item(38, 4) "0.05336"
$ws.Cells.Item(38, 5).Value = "  +1.15%  "
$ws.Cells.Item(39, 2).Value = "MXToken"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Cells.Item(39, 4) "2.999"
$ws.Cells.Item(39, 5).Value = "  +0.94%  "
$ws.Cells.Item(40, 2).Value = "VeChain"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Cells.Item(40, 4) "0.01954"
$ws.Cells.Item(40, 5).Value = "  +0.27%  "
Set-TextValue $ws.Cells.Item(41, 4) "7.344"
$ws.Cells.Item(41, 5).Value = "  +1.46%  "
Set-TextValue $ws.Cells.Item(42, 4) "2.398"
$ws.Cells.Item(42, 5).Value = "  +5.31%  "
Set-TextValue $ws.Cells.Item(43, 4) "0.5331"
$ws.Cells.Item(43, 5).Value = "  +0.55%  "
$ws.Cells.Item(44, 5).Value = "  +0.39%  "
Set-TextValue $ws.Cells.Item(45, 4) "8.512"
$ws.Cells.Item(45, 5).Value = "  +1.06%  "
Set-TextValue $ws.Cells.Item(46, 4) "0.4916"
$ws.Cells.Item(46, 5).Value = "  +0.84%  "
Set-TextValue $ws.Cells.Item(47, 4) "10.59"
$ws.Cells.Item(47, 5).Value = "  +2.08%  "
Set-TextValue $ws.Cells.Item(48, 4) "105.17"
$ws.Cells.Item(48, 5).Value = "  +1.94%  "
$ws.Cells.Item(49, 5).Value = "  +1.02%  "
Set-TextValue $ws.Cells.Item(50, 4) "1.678"
$ws.Cells.Item(50, 5).Value = "  +1.19%  "
Set-TextValue $ws.Cells.Item(51, 4) "0.06320"
$ws.Cells.Item(51, 5).Value = "  +0.45%  "
